# Update quizvragen via Admin
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DC")

# 1. Update choices text in E2 (add leading spaces to 2nd/3rd options)
$ws.Range("E2").Value = "['I = U/R', '     U = I*R', '     R = U/I']"

# 2. Clear the image_url in L2 (question's image removed)
$ws.Range("L2").Value = ""
$ws.Range("L2").Style = "Normal"

# 3. Add a brand-new question row (row 6)
$ws.Range("A6:L6").Style = "Normal"
$ws.Range("B6").Value = "input"
$ws.Range("D6").Value = "4 Wat is er aan de gang?"
$ws.Range("F6").Value = "Klote"
$ws.Range("L6").Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/DC_new_1763126859.png"
